$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.884.11'
$ws.Range("E2").Value = '  -0.83%  '

$ws.Range("D3").Value = '2.294.71'
$ws.Range("E3").Value = '  -1.12%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.78'
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.50'
$ws.Range("E6").Value = '  -3.87%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.508'
$ws.Range("E7").Value = '  -0.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.36'
$ws.Range("E10").Value = '  -4.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.15'
$ws.Range("E12").Value = '  -4.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.83'
$ws.Range("E14").Value = '  +7.57%  '

$ws.Range("E15").Value = '  -0.24%  '

$ws.Range("D16").Value = '2.651.45'
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("D17").Value = '2.318.78'
$ws.Range("E17").Value = '  +0.43%  '

$ws.Range("E18").Value = '  +0.56%  '

$ws.Range("D19").Value = '42.771.73'
$ws.Range("E19").Value = '  -0.85%  '

$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.52'
$ws.Range("E21").Value = '  -2.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.02'
$ws.Range("E22").Value = '  -1.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.42'
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.73'
$ws.Range("E24").Value = '  -0.79%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -2.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.28'
$ws.Range("E28").Value = '  -2.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.04'
$ws.Range("E29").Value = '  +1.20%  '

$ws.Range("E30").Value = '  -5.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.61'
$ws.Range("E31").Value = '  -2.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.11'
$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.69'
$ws.Range("E34").Value = '  +4.96%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.92'
$ws.Range("E35").Value = '  -2.59%  '

$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '16.92'
$ws.Range("E36").Value = '  +0.77%  '

$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.38'
$ws.Range("E37").Value = '  -1.79%  '

$ws.Range("E39").Value = '  -2.80%  '

$ws.Range("E40").Value = '  -1.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").Value = '  -3.73%  '

$ws.Range("E42").Value = '  -1.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.34'
$ws.Range("E43").Value = '  -2.84%  '

$ws.Range("D44").Value = '1.992.45'
$ws.Range("E44").Value = '  +0.79%  '

$ws.Range("E45").Value = '  -1.33%  '

$ws.Range("E46").Value = '  +0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.53'
$ws.Range("E47").Value = '  -5.14%  '

$ws.Range("E48").Value = '  -3.15%  '

$ws.Range("D49").Value = '2.520.54'
$ws.Range("E49").Value = '  -1.14%  '

$ws.Range("E50").Value = '  -3.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.56'
$ws.Range("E51").Value = '  -5.87%  '
